$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1800193333333333
$ws.Range("H2").Value = 0.5400579999999999
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 26.22895092106355
$ws.Range("R2").Value = 236.060558289572
$ws.Range("S2").Value = 0.2865937750105843
$ws.Range("T2").Value = 0.2865937750105843

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1800193333333333
$ws.Range("H3").Value = 0.5400579999999999
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 30.38720994091822
$ws.Range("R3").Value = 273.484889468264
$ws.Range("S3").Value = 0.3320294904365841
$ws.Range("T3").Value = 0.3320294904365841

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1800193333333333
$ws.Range("H4").Value = 0.5400579999999999
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 23.06518494565689
$ws.Range("R4").Value = 207.586664510912
$ws.Range("S4").Value = 0.2520245069956105
$ws.Range("T4").Value = 0.2520245069956105

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1800193333333333
$ws.Range("H5").Value = 0.5400579999999999
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 11.83826560086066
$ws.Range("R5").Value = 106.544390407746
$ws.Range("S5").Value = 0.1293522275572212
$ws.Range("T5").Value = 0.1293522275572212

$wb.Save()
